$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new feature value for the "Строительство" (Construction) category:
# "Введ. кварт. - consnewapt (шт. на 1000 чел.) (id8215002)"
$ws.Range("C12").Value = "Введ. кварт. - consnewapt (шт. на 1000 чел.) (id8215002)"

# Extend the orange category-label fill from B4 down into B5, and likewise for B12,
# so that the "признаки" label cells visually continue onto the row below,
# while keeping left alignment (not centered like B4).
$ws.Range("B5").Interior.ThemeColor = 10
$ws.Range("B5").Interior.TintAndShade = 0.59999389629810485

$ws.Range("B12").Interior.ThemeColor = 10
$ws.Range("B12").Interior.TintAndShade = 0.59999389629810485

# Update the last active selection to reflect where the user finished editing.
$ws.Range("D20").Select()
